$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The edit: the "_GoBack" bookmark moves from the end of paragraph 1 (the
# title) to the end of paragraph 2 (the "1. Introduction" heading), and the
# heading paragraph gains a new "5" appended to its text ("1. Introduction"
# -> "1. Introduction5").
# ---------------------------------------------------------------------------

# Step 1: remove the _GoBack bookmark from its current location at the end
# of paragraph 1.
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# Step 2: append "5" to the "1. Introduction" heading (paragraph 2), right
# before its paragraph mark, inheriting the surrounding (bold, 18pt, Times
# New Roman) run formatting.
$p2 = $d.Paragraphs(2)
$insPos = $p2.Range.End - 1
$insR = $d.Range($insPos, $insPos)
$insR.InsertBefore("5")

# Step 3: re-anchor the _GoBack bookmark at the end of paragraph 2 (right
# after the new "5", before the paragraph mark). A collapsed range sitting
# exactly one character before a paragraph's end cannot be used directly as
# the Bookmarks.Add() target, so first add it at the (safe) paragraph
# boundary after inserting a one-character placeholder there, then delete
# the placeholder - the now-existing bookmark tracks the surrounding text
# edit correctly and ends up exactly where required.
$p2 = $d.Paragraphs(2)
$bmTargetPos = $p2.Range.End - 1

$phR = $d.Range($bmTargetPos, $bmTargetPos)
$phR.InsertBefore("Z")

$bmR = $d.Range($bmTargetPos, $bmTargetPos)
$d.Bookmarks.Add("_GoBack", $bmR)

$delR = $d.Range($bmTargetPos, $bmTargetPos + 1)
$delR.Delete()
